# Applies the edits described by the commit diff:
#
#  1. The paragraph that introduces the project-creation steps
#     ("Seguimos los pasos de la linea de comandos ...") gets the
#     "Intense Quote" paragraph style applied, i.e. it gains a
#     <w:pStyle w:val="IntenseQuote"/> in its <w:pPr>.
#
#  2. The "Requisitios" Heading 1 paragraph had its text split across two
#     runs ("R" + "equisitios"). It is normalized back into a single run
#     containing the full text "Requisitios" (no visible text changes).

$d = $word.ActiveDocument

# --- 1. Apply the "Intense Quote" style to the "Seguimos los pasos..." paragraph ---
$scope1 = $d.Content
$found1 = $scope1.Find.Execute("Seguimos los pasos", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $scope1.Paragraphs(1).Style = "Intense Quote"
}

# --- 2. Merge the "R" + "equisitios" runs back into a single "Requisitios" run ---
$scope2 = $d.Content
$scope2.Find.Execute("Requisitios", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Requisitios", 2)
